$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3 each lost their first data column ("Version" header / "1" value),
# so every remaining cell shifts one column to the left (B->A, C->B, D->C, E->D, F->E).
# Copy cell-by-cell (left to right) so we never read from a cell we already overwrote,
# and use Range.Copy (not a re-typed literal) so the original value type/format
# (e.g. the literal text "TRUE") survives instead of being re-interpreted.
$cols = @("A","B","C","D","E","F")
foreach ($row in @(2,3)) {
    for ($i = 0; $i -lt 5; $i++) {
        $destCol = $cols[$i]
        $srcCol = $cols[$i+1]
        $ws.Range($srcCol + $row).Copy($ws.Range($destCol + $row))
    }
    $ws.Range("F" + $row).ClearContents()
}

# New trailing (empty) cells appear at I2/J2 and I3/J3 in the edited sheet.
$ws.Cells.Item(2,9).Value = $null
$ws.Cells.Item(2,10).Value = $null
$ws.Cells.Item(3,9).Value = $null
$ws.Cells.Item(3,10).Value = $null

# Selection moved from A4 to A2.
$ws.Range("A2").Select()
